# Add 2022-Q3 data
# 1) Update the "总计" (summary) sheet: insert a new row for 2022-Q3 right
#    after the header, shifting the existing quarter rows down by one.
# 2) Add a new "2022-Q3" worksheet (positioned before "2021-Q4") containing
#    the underlying fund-holding detail for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: "总计" sheet - insert the new 2022-Q3 summary row
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0

# Restore the index-column style (bold / bordered / centered) that the
# other "A" cells in this column use, by copying it from the row below.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A is a plain 0-based row counter (not quarter-linked data), so
# it must stay sequential 0,1,2,3 top-to-bottom rather than shifting down
# together with the B/C/D quarter data that the row-insert just moved.
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# ---------------------------------------------------------------------
# Part 2: new "2022-Q3" worksheet with the fund-holding detail table
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")

# Duplicate the "2021-Q4" sheet (this keeps all formatting/styles intact)
# and place the copy immediately before it, then rename + re-point it at
# the new quarter's data.
$q4.Copy($q4)
$q3new = $wb.Worksheets.Item("2021-Q4 (2)")
$q3new.Name = "2022-Q3"

# The source sheet had 3 data rows; 2022-Q3 only needs 2, so drop the
# extra row (the 3rd data row, sheet row 4).
$q3new.Rows.Item(4).Delete()

# Row 2 (fund 005105)
$q3new.Range("A2").Value = 0
$q3new.Range("B2").Value = "'005105"
$q3new.Range("C2").Value = "富荣福康混合C"
$q3new.Range("D2").Value = "'0.07"
$q3new.Range("E2").Value = "'91.00"
$q3new.Range("F2").Value = "'3.04"
$q3new.Range("G2").Value = "'0.0021"
$q3new.Range("H2").Value = 10

# Row 3 (fund 005104)
$q3new.Range("A3").Value = 1
$q3new.Range("B3").Value = "'005104"
$q3new.Range("C3").Value = "富荣福康混合A"
$q3new.Range("D3").Value = "'0.03"
$q3new.Range("E3").Value = "'91.00"
$q3new.Range("F3").Value = "'3.04"
$q3new.Range("G3").Value = "'0.0009"
$q3new.Range("H3").Value = 10

# The apostrophe-forced-text cells above pick up a "quote prefix" style;
# reset their style back to Normal so they match the plain data cells in
# the other quarter sheets (text type is preserved, formatting is not).
$q3new.Range("B2:G3").Style = "Normal"

Write-Output "2022-Q3 sheet added and 总计 sheet updated"
